$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-16 Sunday" "2025-02-17 Monday"
Replace-Text "517×3=1551" "955×2=1910"
Replace-Text "626×8=5008" "977×4=3908"
Replace-Text "648×8=5184" "820×7=5740"
Replace-Text "888×9=7992" "553×9=4977"
Replace-Text "197×7=1379" "827×4=3308"
Replace-Text "333×6=1998" "482×7=3374"
Replace-Text "803×6=4818" "470×4=1880"
Replace-Text "722×5=3610" "375×7=2625"
Replace-Text "779×3=2337" "591×9=5319"
Replace-Text "498×5=2490" "430×3=1290"
Replace-Text "161×5=805" "842×2=1684"
Replace-Text "680×8=5440" "975×7=6825"
Replace-Text "612×8=4896" "662×3=1986"
Replace-Text "237×5=1185" "579×5=2895"
Replace-Text "356×6=2136" "438×4=1752"
Replace-Text "252×7=1764" "914×3=2742"
Replace-Text "164×4=656" "353×7=2471"
Replace-Text "503×4=2012" "456×6=2736"
Replace-Text "353×9=3177" "626×6=3756"
Replace-Text "884×3=2652" "250×8=2000"
Replace-Text "233×6=1398" "927×4=3708"
Replace-Text "728×3=2184" "266×3=798"
Replace-Text "125×3=375" "492×3=1476"
Replace-Text "782×4=3128" "877×5=4385"
Replace-Text "631×3=1893" "457×7=3199"
